# Generate Report for Handoff
# Refresh the "latest handoff" timestamps that are reported in the
# localization-status workbook, simulating a freshly regenerated report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: column D is "Latest Handoff Date". All rows that were
# still showing the stale "2016-15-11 08:15:12" / "2016-15-11 08:15:48"
# timestamps now get the new, single, regenerated timestamp.
$overviewOld1 = "2016-15-11 08:15:12"
$overviewOld2 = "2016-15-11 08:15:48"
$overviewNew  = "2016-16-11 08:16:07"

$overviewUsed = $wsOverview.UsedRange
for ($r = 1; $r -le $overviewUsed.Rows.Count; $r++) {
    $cell = $wsOverview.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -eq $overviewOld1 -or $val -eq $overviewOld2) {
        $cell.Value2 = $overviewNew
    }
}

# zh-cn sheet: column E is "Latest Handoff Datetime".
$zhOld1 = "2016-03-11 08:15:09"
$zhOld2 = "2016-03-11 08:15:42"
$zhNew  = "2016-03-11 08:16:02"

$zhUsed = $wsZhCn.UsedRange
for ($r = 1; $r -le $zhUsed.Rows.Count; $r++) {
    $cell = $wsZhCn.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -eq $zhOld1 -or $val -eq $zhOld2) {
        $cell.Value2 = $zhNew
    }
}

# de-de sheet: column E is "Latest Handoff Datetime".
$deOld1 = "2016-03-11 08:15:12"
$deOld2 = "2016-03-11 08:15:48"
$deNew  = "2016-03-11 08:16:07"

$deUsed = $wsDeDe.UsedRange
for ($r = 1; $r -le $deUsed.Rows.Count; $r++) {
    $cell = $wsDeDe.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -eq $deOld1 -or $val -eq $deOld2) {
        $cell.Value2 = $deNew
    }
}
